# Atividade 3 Questao 4 - Solver
# Update the solver inputs (B2:D2) with a feasible solution and extend the
# solver_lhs3 / solver_rhs3 named ranges to cover rows 9:11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Extend the solver constraint ranges (solver_lhs3 / solver_rhs3) to B9:B11 / D9:D11
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Planilha1!solver_lhs3") {
        $n.RefersTo = "=Planilha1!`$B`$9:`$B`$11"
    }
    elseif ($n.Name -eq "Planilha1!solver_rhs3") {
        $n.RefersTo = "=Planilha1!`$D`$9:`$D`$11"
    }
}

# Set the decision variables found by the Solver
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 11
$ws.Range("D2").Value = 56

# Match the selection saved in the workbook
[void]$ws.Range("E3").Select()
